$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New error code rows 38-41: code only, no message/level
$ws.Range("A38").Value = 10036
$ws.Range("A39").Value = 10037
$ws.Range("A40").Value = 10038
$ws.Range("A41").Value = 10039

# Row 42: new physical address record created successfully message
$ws.Range("A42").Value = 10040
$ws.Range("B42").Value = "message_10040_physical_address_record_created_successfully"
$ws.Range("D42").Value = "Success"

# Update the selection to match the edited cell
$ws.Range("B42").Select()
